$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the cryptos list: latest prices and 1h volume deltas, plus two
# coin rows that swapped rank position (ImmutableX/TrustWalletToken and
# Aave/PaxDollar).

$ws.Range('D2').Value = '27.912.73'
$ws.Range('E2').Value = '  +1.40%  '
$ws.Range('D3').Value = '1.643.66'
$ws.Range('E3').Value = '  +1.31%  '
$ws.Range('D5').Value = '''213.64'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +1.04%  '
$ws.Range('D6').Value = '''0.525'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  -0.27%  '
$ws.Range('D8').Value = '''23.45'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  +1.32%  '
$ws.Range('D9').Value = '''0.264'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  +0.35%  '
$ws.Range('E10').Value = '  +0.67%  '
$ws.Range('D11').Value = '''0.0877'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  -1.26%  '
$ws.Range('D12').Value = '1.876.66'
$ws.Range('E12').Value = '  +1.34%  '
$ws.Range('D13').Value = '1.647.20'
$ws.Range('E13').Value = '  +1.52%  '
$ws.Range('D14').Value = '''0.572'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  +4.23%  '
$ws.Range('E15').Value = '  +0.53%  '
$ws.Range('D16').Value = '''65.88'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  +1.00%  '
$ws.Range('D17').Value = '27.895.59'
$ws.Range('E17').Value = '  +1.48%  '
$ws.Range('D18').Value = '''230.37'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  -0.77%  '
$ws.Range('D19').Value = '0.0₃0724'
$ws.Range('E19').Value = '  +0.63%  '
$ws.Range('E20').Value = '  +1.09%  '
$ws.Range('E21').Value = '  +0.10%  '
$ws.Range('D22').Value = '''10.80'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  +5.91%  '
$ws.Range('D23').Value = '''4.40'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  +1.70%  '
$ws.Range('D24').Value = '''2.13'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  +2.71%  '
$ws.Range('D25').Value = '''152.72'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  +1.93%  '
$ws.Range('E26').Value = '  +0.66%  '
$ws.Range('E27').Value = '  +0.94%  '
$ws.Range('D28').Value = '''15.68'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  +0.94%  '
$ws.Range('E29').Value = '  +0.04%  '
$ws.Range('E30').Value = '  +1.09%  '
$ws.Range('E31').Value = '  +0.60%  '
$ws.Range('D32').Value = '''3.34'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  +1.92%  '
$ws.Range('D33').Value = '1.435.68'
$ws.Range('E33').Value = '  -3.03%  '
$ws.Range('E34').Value = '  -0.20%  '
$ws.Range('E35').Value = '  +1.50%  '
$ws.Range('E36').Value = '  +0.05%  '
$ws.Range('D37').Value = '''0.885'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  +1.61%  '
$ws.Range('E38').Value = '  +0.86%  '
$ws.Range('B39').Value = 'TrustWalletToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D39').Value = '''0.926'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  -2.85%  '
$ws.Range('B40').Value = 'ImmutableX'
$ws.Range('C40').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D40').Value = '''0.559'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  +1.07%  '
$ws.Range('E41').Value = '  +1.56%  '
$ws.Range('B42').Value = 'PaxDollar'
$ws.Range('C42').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D42').Value = '''1.00'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  +0.00%  '
$ws.Range('B43').Value = 'Aave'
$ws.Range('C43').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D43').Value = '''68.60'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  +1.39%  '
$ws.Range('D44').Value = '''2.47'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  -0.06%  '
$ws.Range('D45').Value = '''5.43'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  +4.05%  '
$ws.Range('E46').Value = '  +4.25%  '
$ws.Range('E47').Value = '  +0.35%  '
$ws.Range('D48').Value = '1.785.37'
$ws.Range('E48').Value = '  +1.27%  '
$ws.Range('D49').Value = '''89.18'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  +2.23%  '
$ws.Range('E50').Value = '  -0.04%  '
$ws.Range('D51').Value = '''7.71'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  +0.16%  '

Write-Host "Applied cryptos update"
